$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the date-serial values in column A (rows 2-39) with quarterly
# text labels like "1987Q4" .. "2024Q4", and make them look like the
# header cell (same style/format, no more custom date-time number format).

$startYear = 1987
for ($i = 0; $i -lt 38; $i++) {
    $row = 2 + $i
    $year = $startYear + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "$($year)Q4"
}

# Make column A rows 2-39 match the style of the header cell A1 (which
# uses the shared "centered / bordered" style without any custom
# date-time number format) by copying A1's format over them — this
# reuses the existing style instead of allocating a new one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$wb.Save()
